$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# TC2's step 2 (row 20) currently holds the "atribuir/desatribuir" text,
# TC3's step 2 (row 28) currently holds the "realizar a liquidação" text.
# Swap them so TC2 gets the "liquidação" step and TC3 gets the
# "atribuir/desatribuir" step.

$tc2Step   = $ws.Range("B20").Value2
$tc2Result = $ws.Range("D20").Value2
$tc3Step   = $ws.Range("B28").Value2
$tc3Result = $ws.Range("D28").Value2

$ws.Range("B20").Value = $tc3Step
$ws.Range("D20").Value = $tc3Result
$ws.Range("B28").Value = $tc2Step
$ws.Range("D28").Value = $tc2Result
